$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column S (02-jul) with header + values
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (R1) onto the new header
# cell (S1) so the new column matches the existing bold/bordered header style.
$wsSpot.Range("R1").Copy()
$wsSpot.Range("S1").PasteSpecial(-4122)  # xlPasteFormats
$wsSpot.Range("S1").Value = "02-jul"

$spotValues = @{
    2  = 119.29
    3  = 102.71
    4  = 92.54
    5  = 92.08
    6  = 85.08
    7  = 93.47
    8  = 103.68
    9  = 107.8
    10 = 104.52
    11 = 92.73
    12 = 81.79
    13 = 76.39
    14 = 70.64
    15 = 50.29
    16 = 54.02
    17 = 71.57
    18 = 79.98
    19 = 87.61
    20 = 87.68
    21 = 119.03
    22 = 119.8
    23 = 126.16
    24 = 125.1
    25 = 108.51
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 19).Value = $spotValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 16 (2025-06-30, 31.325)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date-looking text to be stored as plain text (matching the rest
# of column A) instead of letting Excel auto-convert it to a date serial.
$wsGaz.Range("A16").NumberFormat = "@"
$wsGaz.Range("A16").Value = "2025-06-30"
# Re-apply the plain (unstyled) formatting used by the rest of column A.
$wsGaz.Range("A2").Copy()
$wsGaz.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$wsGaz.Range("B16").Value = 31.325

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 16 (2025-06-30, 68)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A16").NumberFormat = "@"
$wsCo2.Range("A16").Value = "2025-06-30"
$wsCo2.Range("A2").Copy()
$wsCo2.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$wsCo2.Range("B16").Value = 68
